# Update the "Metadata" worksheet of the HIPAA relationship value-set
# workbook to reflect the new IG publication metadata (version bump,
# regenerated date, publisher/jurisdiction in place of the old
# contact-detail rows) and drop the now-removed trailing "Immutable" row
# duplication (the sheet shrinks from 15 to 14 data rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Version: 5.0.0 -> 6.0.0
$ws.Range("B3").Value = "6.0.0"

# Date: regenerated publication timestamp
$ws.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher row gains a value; "Contact" / "No display for ContactDetail"
# rows are replaced by "Jurisdiction" / "United States of America", which
# shifts the remaining rows up by one and drops the final duplicate row.
$ws.Range("B9").Value = "Alvearie Team"

$ws.Range("A10").Value = "Jurisdiction"
$ws.Range("B10").Value = "United States of America"

$ws.Range("A11").Value = "Description"
$ws.Range("B11").Value = "HIPAA standard value set for the relationship of the member to the contract holder"

$ws.Range("A12").Value = "Purpose"
$ws.Range("B12").Value = ""

$ws.Range("A13").Value = "Copyright"
$ws.Range("B13").Value = ""

$ws.Range("A14").Value = "Immutable"
$ws.Range("B14").Value = "BooleanType[null]"

# The old row 15 (duplicate "Immutable"/"BooleanType[null]") no longer
# exists; delete the now-unused trailing row so the sheet's used range
# shrinks to A1:B14.
$ws.Rows.Item(15).Delete()
